# Applies the "added Period table" edit described by the diff:
#   - Paragraph "Weekday (weekdayID, name)" becomes
#     "Period (periodID, starttime, endtime)".
#   - Paragraph "Schedule (scheduleID, name)" becomes
#     "Weekday (weekdayID, name)".
#   - A new paragraph "Schedule (scheduleID, name)" is inserted right
#     after the (new) Weekday paragraph.
#   - The "subjectsOnWeekday (^subjectID, ^weekdayID)" paragraph gains a
#     third underlined foreign key: ", ^periodID".

$d = $word.ActiveDocument

function New-PackageXml([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1. "Weekday (weekdayID, name)" -> "Period (periodID, starttime, endtime)"
# ---------------------------------------------------------------------
$periodPara = $d.Paragraphs(8)
$periodBody =
    '<w:p>' +
      '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Period (</w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>periodID</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>starttime</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, endtime</w:t></w:r>' +
      '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
      '<w:bookmarkEnd w:id="0"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' +
    '</w:p>'
$periodPara.Range.InsertXML((New-PackageXml $periodBody))

# ---------------------------------------------------------------------
# 2. "Schedule (scheduleID, name)" -> "Weekday (weekdayID, name)"
# ---------------------------------------------------------------------
$weekdayPara = $d.Paragraphs(9)
$weekdayBody =
    '<w:p>' +
      '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Weekday</w:t></w:r>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>weekdayID</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, name)</w:t></w:r>' +
    '</w:p>'
$weekdayPara.Range.InsertXML((New-PackageXml $weekdayBody))

# ---------------------------------------------------------------------
# 3. Insert a brand-new "Schedule (scheduleID, name)" paragraph right
#    after the (now) Weekday paragraph.
# ---------------------------------------------------------------------
$weekdayPara.Range.InsertParagraphAfter()
$schedulePara = $d.Paragraphs(10)
$scheduleBody =
    '<w:p>' +
      '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Schedule</w:t></w:r>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (</w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>scheduleID</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>, name</w:t></w:r>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' +
    '</w:p>'
$schedulePara.Range.InsertXML((New-PackageXml $scheduleBody))

# ---------------------------------------------------------------------
# 4. "subjectsOnWeekday (^subjectID, ^weekdayID)" ->
#    "subjectsOnWeekday (^subjectID, ^weekdayID, ^periodID)"
# ---------------------------------------------------------------------
$subjOnWeekdayPara = $d.Paragraphs(11)
$subjOnWeekdayBody =
    '<w:p>' +
      '<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
      '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>subjectsOnWeekday</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r>' +
      '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>^</w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>subjectID</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' +
      '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>^</w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>weekdayID</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>, ^</w:t></w:r>' +
      '<w:proofErr w:type="spellStart"/>' +
      '<w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>periodID</w:t></w:r>' +
      '<w:proofErr w:type="spellEnd"/>' +
      '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>' +
    '</w:p>'
$subjOnWeekdayPara.Range.InsertXML((New-PackageXml $subjOnWeekdayBody))

Write-Output ("Paragraph count: " + $d.Paragraphs.Count)
foreach ($p in $d.Paragraphs) {
    Write-Output $p.Range.Text
}
